$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("C5").Value = 44
$ws.Range("D5").Value = "Worked on the TODO list and got it finsihed. Going to start on the Bookmark lab. Gonna work on that tomorrow. "

$ws.Range("D5").Select()
